$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append the new ticker to the list
$ws.Range("A30").Value = "BKNG"

# Update the active selection as it was left after the edit
$ws.Range("B32").Select()
